# Bump the "Förändrad" (Changed) date column by one day for every data row.
# Column C holds a date serial number (formatted as YYYY-MM-DD) that was
# incremented from 46075 (2026-02-22) to 46076 (2026-02-23) for all rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row based on column A (Beteckning), data starts at row 2.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = [double]$current + 1
    }
}
